$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.529.59'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '1.913.08'
$ws.Range('E3').Value = '  +0.36%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''0.708'
$ws.Range('E5').Value = '  +9.14%  '
$ws.Range('D6').Value = '''247.05'
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '''40.89'
$ws.Range('E8').Value = '  -2.75%  '
$ws.Range('D9').Value = '''0.357'
$ws.Range('E9').Value = '  +4.29%  '
$ws.Range('D10').Value = '''52.70'
$ws.Range('E10').Value = '  +8.03%  '
$ws.Range('D11').Value = '''0.0731'
$ws.Range('E11').Value = '  +3.08%  '
$ws.Range('D12').Value = '''0.0990'
$ws.Range('E12').Value = '  -0.95%  '
$ws.Range('D13').Value = '2.190.73'
$ws.Range('D14').Value = '''12.64'
$ws.Range('E14').Value = '  +1.73%  '
$ws.Range('D15').Value = '''0.716'
$ws.Range('E15').Value = '  +2.47%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.919.77'
$ws.Range('E16').Value = '  -0.32%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = '''4.91'
$ws.Range('E17').Value = '  +1.58%  '
$ws.Range('D18').Value = '35.534.07'
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('D19').Value = '''73.32'
$ws.Range('E19').Value = '  +1.65%  '
$ws.Range('E20').Value = '  -0.66%  '
$ws.Range('D21').Value = '''13.15'
$ws.Range('E21').Value = '  +3.75%  '
$ws.Range('D22').Value = '''242.70'
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('E23').Value = '  +4.09%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('D26').Value = '''2.30'
$ws.Range('E26').Value = '  +2.38%  '
$ws.Range('D27').Value = '''168.56'
$ws.Range('E27').Value = '  -1.84%  '
$ws.Range('D28').Value = '''8.66'
$ws.Range('E28').Value = '  +1.45%  '
$ws.Range('D29').Value = '''18.77'
$ws.Range('E29').Value = '  +4.04%  '
$ws.Range('E30').Value = '  +4.47%  '
$ws.Range('E32').Value = '  +2.20%  '
$ws.Range('E33').Value = '  +0.59%  '
$ws.Range('E34').Value = '  +0.34%  '
$ws.Range('E35').Value = '  +6.76%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('E37').Value = '  -4.82%  '
$ws.Range('E38').Value = '  +11.48%  '
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('D40').Value = '''17.31'
$ws.Range('E40').Value = '  +10.28%  '
$ws.Range('D41').Value = '''97.90'
$ws.Range('E41').Value = '  +5.71%  '
$ws.Range('E42').Value = '  +1.19%  '
$ws.Range('E43').Value = '  +2.63%  '
$ws.Range('D44').Value = '''0.0650'
$ws.Range('E44').Value = '  +1.42%  '
$ws.Range('D45').Value = '1.354.88'
$ws.Range('E45').Value = '  +0.48%  '
$ws.Range('E46').Value = '  +2.39%  '
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('B48').Value = 'MXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D48').Value = '''2.80'
$ws.Range('E48').Value = '  +1.23%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').Value = '''46.04'
$ws.Range('E49').Value = '  -6.18%  '
$ws.Range('D50').Value = '''12.27'
$ws.Range('E50').Value = '  -3.25%  '
$ws.Range('E51').Value = '  -0.68%  '
